# Price-scrape update: append the 2026-02-07 row to the price-history sheet.
# Source row (35) holds: Date | Price | Discount | Incredible, all stored as
# plain text (shared strings), e.g. "2025-11-21" | "17799000" | "0" | "0".
# The new row 36 continues that pattern with "2026-02-07" | "26349000" | "0" | "0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 36
$rng = $ws.Range("A" + $newRow + ":D" + $newRow)

# Force the new cells to Text so Excel stores the date/number-looking
# strings verbatim instead of auto-converting them to a date serial /
# numeric value (matching how the rest of the column is stored).
$rng.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "26349000"
$ws.Cells.Item($newRow, 3).Value = "0"
$ws.Cells.Item($newRow, 4).Value = "0"

# Drop the temporary Text number format again so the new cells end up with
# the workbook's default (unstyled) formatting, same as every other row.
$rng.ClearFormats()
